$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The "Is in use" column (J) used to hold text values "no"/"yes" (shared strings).
# Data imports now expect numeric booleans instead: 0 = no, 1 = yes.
$ws.Range("J2").Value = 0
$ws.Range("J3").Value = 1

# Reflect the last active selection being on column J (row 4) as last saved by Excel.
$ws.Range("J4").Select()
